$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = @(-19.10121911200491, 2.148082449484106, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491)
for ($c = 0; $c -lt $row2.Length; $c++) {
    $ws.Cells.Item(2, $c + 2).Value = $row2[$c]
}

$row3 = @(-19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, 2.452858648446483, -19.10121911200491, -19.10121911200491)
for ($c = 0; $c -lt $row3.Length; $c++) {
    $ws.Cells.Item(3, $c + 2).Value = $row3[$c]
}

$row4 = @(-19.10121911200491, 2.096230975789803, 2.692209961642769, -19.10121911200491, 2.593412958216263, -19.10121911200491, 1.803572043601532, -19.10121911200491, 2.406523758597115, -19.10121911200491)
for ($c = 0; $c -lt $row4.Length; $c++) {
    $ws.Cells.Item(4, $c + 2).Value = $row4[$c]
}

$row5 = @(-19.10121911200491, 1.182408495485197, -19.10121911200491, -19.10121911200491, -19.10121911200491, 2.109845119594896, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491)
for ($c = 0; $c -lt $row5.Length; $c++) {
    $ws.Cells.Item(5, $c + 2).Value = $row5[$c]
}

$row6 = @(-19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491)
for ($c = 0; $c -lt $row6.Length; $c++) {
    $ws.Cells.Item(6, $c + 2).Value = $row6[$c]
}

$row7 = @(2.969747241494152, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491)
for ($c = 0; $c -lt $row7.Length; $c++) {
    $ws.Cells.Item(7, $c + 2).Value = $row7[$c]
}

$row8 = @(-19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491)
for ($c = 0; $c -lt $row8.Length; $c++) {
    $ws.Cells.Item(8, $c + 2).Value = $row8[$c]
}

$row9 = @(3.60478042374679, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491)
for ($c = 0; $c -lt $row9.Length; $c++) {
    $ws.Cells.Item(9, $c + 2).Value = $row9[$c]
}

$row10 = @(-19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, 1.563305428938405, -19.10121911200491, 2.222546513090658)
for ($c = 0; $c -lt $row10.Length; $c++) {
    $ws.Cells.Item(10, $c + 2).Value = $row10[$c]
}

$row11 = @(-19.10121911200491, -19.10121911200491, -19.10121911200491, 4.321925657870105, -19.10121911200491, 2.553296737039967, -19.10121911200491, -19.10121911200491, -19.10121911200491, 1.316127640993138)
for ($c = 0; $c -lt $row11.Length; $c++) {
    $ws.Cells.Item(11, $c + 2).Value = $row11[$c]
}

$row12 = @(-19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491)
for ($c = 0; $c -lt $row12.Length; $c++) {
    $ws.Cells.Item(12, $c + 2).Value = $row12[$c]
}

$row13 = @(-19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, 2.288028317569583, 1.643250964795983)
for ($c = 0; $c -lt $row13.Length; $c++) {
    $ws.Cells.Item(13, $c + 2).Value = $row13[$c]
}

$row14 = @(-19.10121911200491, -19.10121911200491, 1.614184075505974, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, 2.049292396228972)
for ($c = 0; $c -lt $row14.Length; $c++) {
    $ws.Cells.Item(14, $c + 2).Value = $row14[$c]
}

$row15 = @(-19.10121911200491, -19.10121911200491, 0.2909158993567536, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491)
for ($c = 0; $c -lt $row15.Length; $c++) {
    $ws.Cells.Item(15, $c + 2).Value = $row15[$c]
}

$row16 = @(-19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, 2.296230437378346, -19.10121911200491)
for ($c = 0; $c -lt $row16.Length; $c++) {
    $ws.Cells.Item(16, $c + 2).Value = $row16[$c]
}

$row17 = @(-19.10121911200491, 0.664916679553523, -0.01597229904361772, -19.10121911200491, -19.10121911200491, -19.10121911200491, 0.4887775851738753, 0.8560884359178706, 1.181730155285412, -19.10121911200491)
for ($c = 0; $c -lt $row17.Length; $c++) {
    $ws.Cells.Item(17, $c + 2).Value = $row17[$c]
}

$row18 = @(-19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, -19.10121911200491, 0.4268858361544174, 0.8419896446685495, 1.397214283013171, -19.10121911200491)
for ($c = 0; $c -lt $row18.Length; $c++) {
    $ws.Cells.Item(18, $c + 2).Value = $row18[$c]
}

$row19 = @(-19.10121911200491, -19.10121911200491, 1.999877502719866, -19.10121911200491, -19.10121911200491, -19.10121911200491, 1.87075167403029, 2.060528487064933, -19.10121911200491, -19.10121911200491)
for ($c = 0; $c -lt $row19.Length; $c++) {
    $ws.Cells.Item(19, $c + 2).Value = $row19[$c]
}

$row20 = @(-19.10121911200491, 1.848215081300985, 2.092037320700836, -19.10121911200491, 3.80371478333808, -19.10121911200491, 2.165741813025353, 1.924342917764891, -19.10121911200491, 2.480296410765245)
for ($c = 0; $c -lt $row20.Length; $c++) {
    $ws.Cells.Item(20, $c + 2).Value = $row20[$c]
}

$row21 = @(-19.10121911200491, 1.939643815061356, -19.10121911200491, -19.10121911200491, -19.10121911200491, 3.29480888432791, 2.489890920197907, -19.10121911200491, -19.10121911200491, -19.10121911200491)
for ($c = 0; $c -lt $row21.Length; $c++) {
    $ws.Cells.Item(21, $c + 2).Value = $row21[$c]
}
